# Auto-update draw results: append the new Pick 4 draw row (2025-12-21)
# to the "Results" sheet, mirroring the existing rows (all 5 columns are
# stored as plain text, matching how SheetJS originally wrote this file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.UsedRange.Rows.Count + 1   # 95 -> 96

# Columns A (date string) and C (numeric-looking string) would otherwise be
# auto-coerced into a real date/number by Excel's text-to-value parsing, so
# force them to Text format before assigning, then drop the leftover
# number-format so no extra style lingers on the cell (matches the rest of
# the sheet, which carries no explicit style index).
$cellA = $ws.Cells.Item($newRow, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025-12-21"
$cellA.ClearFormats()

$ws.Cells.Item($newRow, 2).Value = "Pick 4"

$cellC = $ws.Cells.Item($newRow, 3)
$cellC.NumberFormat = "@"
$cellC.Value = "251221"
$cellC.ClearFormats()

$ws.Cells.Item($newRow, 4).Value = "7-8-4-4"

# Column E (ISO timestamp) is not numeric-looking enough to be auto-parsed,
# so it can be assigned directly and stays text.
$ws.Cells.Item($newRow, 5).Value = "2025-12-21T21:38:21.323+04:00"
